# Auto-generated PowerShell COM-interop script implementing the commit
# "fix: contrucción grafico de resultados" against TOPSIS.xlsx.
#
# Summary of the edit:
#  - info!C2: date string corrected 07/11/23 -> 07/17/23
#  - alternative_info: a handful of cells get last-digit floating point
#    refinements (recomputed from the corrected date / reseeded fuzzy draw)
#  - alternatives_norm: the normalized-criteria table (columns C1.2, C1.4,
#    C2.1, C2.2, C2.3) is recalculated with the corrected inputs
#  - result: the TOPSIS ranking is recalculated (new Evaluation scores +
#    a new order), and a new "Alternatives" column (C) is added that
#    mirrors column A so the results chart can plot Evaluation vs.
#    Alternatives correctly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) info!C2 -- fix the recorded run date. Force text storage (NumberFormat
#    "@") so Excel doesn't reinterpret the "dd/mm/yy"-looking string as a
#    serial date; it was (and must remain) a literal string.
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("info")
$cellC2 = $wsInfo.Range("C2")
$cellC2.NumberFormat = "@"
$cellC2.Value2 = "07/17/23"

# ---------------------------------------------------------------------------
# 2) alternative_info -- tiny floating point refinements on a few cells.
# ---------------------------------------------------------------------------
$wsAltInfo = $wb.Worksheets.Item("alternative_info")
$altInfoData = @{
    "G4" = 640.793947259624
    "L6" = 0.0006455764403994414
    "M6" = 175.0574402764045
    "G7" = 640.793947259624
    "G9" = 640.793947259624
    "G12" = 640.793947259624
    "G15" = 640.793947259624
    "G18" = 640.793947259624
    "G21" = 640.793947259624
    "G22" = 2135.979824198747
    "M23" = 81.04021835940149
}
foreach ($cellRef in $altInfoData.Keys) {
    $wsAltInfo.Range($cellRef).Value2 = $altInfoData[$cellRef]
}

# ---------------------------------------------------------------------------
# 3) alternatives_norm -- recalculated normalized criteria matrix
#    (columns B..J for rows 2..24).
# ---------------------------------------------------------------------------
$wsNorm = $wb.Worksheets.Item("alternatives_norm")
$normCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J")
$normData = @{
    2 = @(0.06014777593876288, 0.001565526213971937, 0.0002198582778971346, 0.3345844034605377, 0.02039745084057954, 0.04964547180033088, 0.1054275896195656, 0, 0.2870341709311114)
    3 = @(0.03025183176213168, 0.00450681788870709, 0.000218878412719932, 0.3336005965116349, 0.02271534298155449, 0.1775780337473374, 0.1370558665054353, 0, 0.2391999263454417)
    4 = @(0.04774139495238834, 0.002668510591997619, 0.0001650286867904367, 0.3020858977494591, 0.155066984231224, 0.1718497100780684, 0.1887153854190225, 0, 0.2511548995647225)
    5 = @(0.2010041332104988, 0.1895130129549755, 0.2107570138903742, 0.2009379528126722, 0.02572860276482192, 0.141298650508634, 0.131784487024457, 0.1458649914978946, 0.216473995861971)
    6 = @(0.05937144644710108, 0.0008302032952881483, 0.0001419353244371834, 0.2886148084468784, 0.2874186254808935, 0.1661213864087995, 0.2403749043326096, 0, 0.2631098727840033)
    7 = @(0.2184936964007554, 0.187674705658266, 0.187120500088516, 0.201573993237332, 0.1580802440144914, 0.1355703268393651, 0.1834440059380442, 0.1458649914978946, 0.2284289690812518)
    8 = @(0.3148093435439518, 0.3122298005355343, 0.2488637489443285, 0.1769707172400298, 0.02874186254808935, 0.1050192672699307, 0.1265131075434787, 0.2917299829957891, 0.1937480653785003)
    9 = @(0.03279342286407274, 0.004139156429365196, 0.0001319108251337204, 0.2822149407065723, 0.1562259303017115, 0.2358159910515717, 0.2045295238619573, 0, 0.2272377772718877)
    10 = @(0.1860561611221832, 0.1909836587923431, 0.2599956661877185, 0.1694494283891154, 0.02688754883530939, 0.2052649314821373, 0.1475986254673919, 0.1458649914978946, 0.1925568735691361)
    11 = @(0.04442347435878548, 0.002300849132655725, 0.0001048553843963725, 0.2665475277634776, 0.288577571551381, 0.2300876673823027, 0.2561890427755445, 0, 0.2391927504911685)
    12 = @(0.2035457243124398, 0.1891453514956336, 0.2249370951640692, 0.1744601978376734, 0.1592391900849789, 0.1995366078128683, 0.199258144380979, 0.1458649914978946, 0.2045118467884169)
    13 = @(0.2998613714556362, 0.3137004463729018, 0.2874402451440581, 0.1523631513069525, 0.02990080861857682, 0.168985548243434, 0.1423272459864136, 0.2917299829957891, 0.1698309430856654)
    14 = @(0.2151757758071526, 0.1873070441989241, 0.2064449256530797, 0.1771270045389468, 0.2915908313346484, 0.1938082841435994, 0.2509176632945662, 0.1458649914978946, 0.2164668200076977)
    15 = @(0.3173509346458929, 0.3118621390761924, 0.2605052235066231, 0.1574463935215631, 0.1622524498682463, 0.163257224574165, 0.1939867649000008, 0.2917299829957891, 0.1817859163049462)
    16 = @(0.0003558875855004722, 0.007448109563442244, 0.00001187234700644527, 0.1257617931209708, 0.02503323512252944, 0.3055105956943439, 0.168684143391305, 0, 0.191365681759772)
    17 = @(0.02947550227046987, 0.003771494970023302, 0.000002040341614585504, 0.2053594998150932, 0.2897365176218685, 0.294053948355806, 0.2720031812184793, 0, 0.2152756281983336)
    18 = @(0.1885977522241242, 0.1906159973330011, 0.281938899969236, 0.1335909638463324, 0.1603981361554664, 0.2635028887863716, 0.2150722828239139, 0.1458649914978946, 0.180594724495582)
    19 = @(0.2849133993673206, 0.3151710922102694, 0.340185080703497, 0.1187177416269694, 0.0310597546890643, 0.2329518292169372, 0.1581413844293484, 0.2917299829957891, 0.1459138207928305)
    20 = @(0.200227803718837, 0.1887776900362917, 0.253474687364663, 0.1409920763189318, 0.2927497774051359, 0.2577745651171027, 0.2667318017375011, 0.1458649914978946, 0.1925496977148629)
    21 = @(0.3024029625575773, 0.3133327849135599, 0.3030918200331355, 0.1283013461325413, 0.1634113959387338, 0.2272235055476683, 0.2098009033429356, 0.2917299829957891, 0.1578687940121113)
    22 = @(0.05859511695543927, 0.000094880376604359798, 0.0000019567386732844978, 0.2060363304447788, 0.5544398001212074, 0.2825973010172681, 0.3753222190456537, 0, 0.2391855746368952)
    23 = @(0.31403301405229, 0.3114944776168505, 0.2826207892427517, 0.1336099000497293, 0.2957630371884033, 0.2214951818783993, 0.2614604222565228, 0.2917299829957891, 0.1698237672313921)
    24 = @(0.3676786146310353, 0.4021712540912516, 0.3407803307405586, 0.1187054156359452, 0.03708627425559916, 0.1603930627395305, 0.1475986254673919, 0.5834599659915782, 0.100461959825889)
}
foreach ($r in $normData.Keys) {
    $rowVals = $normData[$r]
    for ($i = 0; $i -lt $normCols.Length; $i++) {
        $wsNorm.Cells.Item([int]$r, $i + 2).Value2 = $rowVals[$i]
    }
}

# ---------------------------------------------------------------------------
# 4) result -- recalculated ranking (new Evaluation values + new row
#    order) and a new "Alternatives" column C mirroring column A, used by
#    the results chart (this is the actual fix for the broken chart
#    referenced in the commit message).
# ---------------------------------------------------------------------------
$wsResult = $wb.Worksheets.Item("result")

$wsResult.Range("C1").Value2 = "Alternatives"
$wsResult.Range("C1").Style = $wsResult.Range("B1").Style

$resultData = @{
    2 = @(13, 0.651095603765667, 13)
    3 = @(21, 0.645652888279656, 21)
    4 = @(6, 0.6443595015312608, 6)
    5 = @(22, 0.6189634769465667, 22)
    6 = @(19, 0.6158451352156784, 19)
    7 = @(11, 0.6121330089224496, 11)
    8 = @(17, 0.5766067563698781, 17)
    9 = @(5, 0.5614105739375657, 5)
    10 = @(12, 0.5569113818910938, 12)
    11 = @(10, 0.5214971392676947, 10)
    12 = @(3, 0.516261064615348, 3)
    13 = @(18, 0.5155967969671299, 18)
    14 = @(16, 0.4802415859157514, 16)
    15 = @(8, 0.475886143715695, 8)
    16 = @(0, 0.4085771840981908, 0)
    17 = @(20, 0.4073584863335043, 20)
    18 = @(4, 0.4012702982355726, 4)
    19 = @(2, 0.3779607510006773, 2)
    20 = @(9, 0.3743831274074745, 9)
    21 = @(15, 0.3530585791620973, 15)
    22 = @(1, 0.3527805223648723, 1)
    23 = @(7, 0.3513393441577106, 7)
    24 = @(14, 0.3217336849476338, 14)
}
foreach ($r in $resultData.Keys) {
    $vals = $resultData[$r]
    $rowIndex = [int]$r
    $wsResult.Cells.Item($rowIndex, 1).Value2 = $vals[0]
    $wsResult.Cells.Item($rowIndex, 2).Value2 = $vals[1]
    $wsResult.Cells.Item($rowIndex, 3).Value2 = $vals[2]
}

Write-Host "Edit applied: info.C2, alternative_info (11 cells), alternatives_norm (23x9), result (23x3 + new column C)"
